# Applies the "frequency parameter optional" data refresh to the
# option-valuation workbook: new Black-Scholes inputs/output, a refreshed
# volatility peer set (now 23 tickers instead of 3), and a refreshed
# risk-free-rate curve.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write literal TEXT (not an auto-parsed Excel date/number) into
# a cell without disturbing its existing cell style. A plain
#   $range.Value = "3/31/2022"
# gets auto-coerced to a date serial, and forcing text via a leading
# apostrophe (or NumberFormat = "@") stamps a brand-new "quote prefixed"
# style onto the cell. Instead, stage the text on a scratch cell, copy
# it, and paste *values only* into the destination -- PasteSpecial
# values-only carries the already-resolved text across without touching
# number format / style, so the destination keeps its original style.
# ---------------------------------------------------------------------
function Set-TextValue($range, [string]$text) {
    $scratch = $range.Worksheet.Range("ZZ1")
    $scratch.Value = "'" + $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)   # xlPasteValues
    $scratch.Clear()
}

# ---------------------------------------------------------------------
# Sheet 1: "Black Scholes"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Black Scholes")

Set-TextValue $ws1.Range("B2") "3/31/2022"
Set-TextValue $ws1.Range("B3") "3/31/2022"
Set-TextValue $ws1.Range("B4") "3/31/2032"
Set-TextValue $ws1.Range("B5") "3/31/2026"

$ws1.Range("B6").Value = 86.29000000000001
$ws1.Range("B7").Value = 86.29000000000001
$ws1.Range("B8").Value = 7
$ws1.Range("B9").Value = 0.0238
$ws1.Range("B10").Value = 0.4888
$ws1.Range("B11").Value = 45.31

# ---------------------------------------------------------------------
# Sheet 2: "Volatility" -- replace the old 3-ticker table (+ header +
# average) with a refreshed 23-ticker table (+ header + average).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Volatility")

$ws2.Range("B1").Value = "2015-03-31 to 2022-03-31"

$volTickers = @(
    @("BKKT", 173.88),
    @("FRGE", 87.97),
    @("VIRT", 37.77),
    @("HOOD", 92.48999999999999),
    @("COIN", 54.53),
    @("TW", 31.39),
    @("MKTX", 31.08),
    @("AMG", 36.12),
    @("APO", 36.22),
    @("ARES", 35.72),
    @("BX", 33.71),
    @("OWL", 39.56),
    @("BPT.L", 46.92),
    @("CG", 34.74),
    @("EQT.ST", 52.54),
    @("HLNE", 35.07),
    @("KKR", 36.1),
    @("EMG.L", 30.71),
    @("PGHN.SW", 25.04),
    @("PHLL.L", 50.09),
    @("STEP", 45.23),
    @("TPG", 28.43)
)

$row = 2
foreach ($entry in $volTickers) {
    $ws2.Cells.Item($row, 1).Value = $entry[0]
    $ws2.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}

# "Average" label/value used to live on row 5; it now lives on row 24
# (row 2 + 22 tickers).
$ws2.Cells.Item($row, 1).Value = "Average"
$ws2.Cells.Item($row, 2).Value = 48.88

# ---------------------------------------------------------------------
# Sheet 3: "Risk Free Rate"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Risk Free Rate")

Set-TextValue $ws3.Range("B1") "2022-03-31"

$rfr = @(
    0.49, 0.97, 1.46, 1.94, 2.42, 2.4, 2.38, 2.37, 2.35, 2.33,
    2.34, 2.34, 2.35, 2.35, 2.36, 2.37, 2.37, 2.38, 2.38, 2.39,
    2.4, 2.4, 2.41, 2.41, 2.42, 2.43, 2.43, 2.44, 2.44, 2.45
)

$row = 2
foreach ($v in $rfr) {
    $ws3.Cells.Item($row, 2).Value = $v
    $row = $row + 1
}
